$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: add P1 = 14, Q1 = 15, matching the style of the existing header cells (O1 has style index 1:
# bold font, bordered, centered alignment)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Rows 2-25: swap values in columns I, K, M, O and add new columns P, Q
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P: new column, value 2
    $ws.Cells.Item($r, 17).Value = 2  # Q: new column, value 2
}
